$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-07-05 Friday" "2024-07-06 Saturday"
Replace-Text "24×85=" "63×75="
Replace-Text "13×83=" "58×53="
Replace-Text "95×33=" "77×15="
Replace-Text "89×81=" "47×95="
Replace-Text "83×93=" "89×93="
Replace-Text "63×37=" "48×22="
Replace-Text "56×38=" "11×46="
Replace-Text "95×94=" "20×70="
Replace-Text "25×97=" "91×21="
Replace-Text "86×83=" "25×64="
Replace-Text "98×95=" "96×17="
Replace-Text "26×97=" "35×63="
Replace-Text "63×60=" "99×18="
Replace-Text "11×24=" "56×24="
Replace-Text "47×21=" "41×84="
Replace-Text "65×26=" "23×50="
Replace-Text "41×20=" "94×84="
Replace-Text "93×24=" "62×95="
Replace-Text "63×26=" "12×87="
Replace-Text "79×60=" "99×11="
Replace-Text "32×13=" "78×27="
Replace-Text "39×14=" "21×39="
Replace-Text "64×79=" "28×77="
Replace-Text "56×19=" "17×42="
Replace-Text "54×89=" "62×88="
